$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The six numbered labels ("2" .. "7") live inside the group "Group 67"
# (top-level shape 5 on the slide) as GroupItems 4 .. 9.
$grp = $s.Shapes.Item(5)
$items = $grp.GroupItems

$items.Item(4).Left = 580.886653543307
$items.Item(4).Top = 422.9966535433071

$items.Item(5).Left = 616.7768897637795
$items.Item(5).Top = 457.92933070866144

$items.Item(6).Left = 657.7609842519685
$items.Item(6).Top = 492.5874409448819

$items.Item(7).Left = 691.5638188976378
$items.Item(7).Top = 526.6199606299213

$items.Item(8).Left = 731.3668897637796
$items.Item(8).Top = 563.3303543307087

$items.Item(9).Left = 766.4698818897638
$items.Item(9).Top = 598.9621653543308
